$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 so existing rows 6-17 shift down to 7-18,
# keeping the same per-row formatting (style s="2" on column D, etc.)
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(6, 3).Value = "Los Lagos"
$ws.Cells.Item(6, 4).Value = 44487
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Value = 300000000
$ws.Cells.Item(6, 7).Value = "Espárragos"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 120
$ws.Cells.Item(6, 11).Value = 1800
$ws.Cells.Item(6, 12).Value = 1800
$ws.Cells.Item(6, 13).Value = 1800
$ws.Cells.Item(6, 14).Value = "$/kilo"
$ws.Cells.Item(6, 15).Value = "Provincia de Linares"
$ws.Cells.Item(6, 16).Value = 1800
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
